$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths for A, B, C (closest achievable to 20.42578125 / 34 / 37.140625 chars) ---
$ws.Range("A1").ColumnWidth = 19.59
$ws.Range("B1").ColumnWidth = 33.17
$ws.Range("C1").ColumnWidth = 36.25

# --- Header row styling: copy bold style (s=1) from A1 to new header cells E1:J1 ---
$ws.Range("A1").Copy()
$ws.Range("E1:J1").PasteSpecial(-4122)

# --- Propagate the "no theme color" style (s=2) used in column D to new rows 10:13 ---
# (row 9 is intentionally left with the default style, matching the source data)
$ws.Range("D2").Copy()
$ws.Range("D10:D13").PasteSpecial(-4122)

# --- Populate all cell values (rows 1-13, cols A-J) ---
# Row 1
$ws.Cells.Item(1,1).Value = "Project "
$ws.Cells.Item(1,2).Value = "Source layer "
$ws.Cells.Item(1,3).Value = "Tileset Name "
$ws.Cells.Item(1,4).Value = "URL "
$ws.Cells.Item(1,5).Value = "Category "
$ws.Cells.Item(1,6).Value = "Type "
$ws.Cells.Item(1,7).Value = "Status "
$ws.Cells.Item(1,8).Value = "Color "
$ws.Cells.Item(1,9).Value = "Shape "
$ws.Cells.Item(1,10).Value = "Shape-Fill"

# Row 2
$ws.Cells.Item(2,1).Value = "Salt Marsh Plant Survey "
$ws.Cells.Item(2,2).Value = "Vegetation_Survey_Plots_LHGSM"
$ws.Cells.Item(2,3).Value = "Vegetation_Survey_Plots_LHGSM"
$ws.Cells.Item(2,4).Value = "jmkhoch.cjsg8if770gxy2rpgqfjufj5e-6ofw0"
$ws.Cells.Item(2,5).Value = "Vegetation "
$ws.Cells.Item(2,6).Value = "Monitoring project "
$ws.Cells.Item(2,7).Value = "Ongoing "
$ws.Cells.Item(2,8).Value = "Green"
$ws.Cells.Item(2,9).Value = "Circle "
$ws.Cells.Item(2,10).Value = "Fill "

# Row 3
$ws.Cells.Item(3,1).Value = "Harbor Herons "
$ws.Cells.Item(3,2).Value = "Wading_Bird_Count_Survey_Points_"
$ws.Cells.Item(3,3).Value = "Wading_Bird_Count_Survey_Points_"
$ws.Cells.Item(3,4).Value = "jmkhoch.cjsg8hlb90fdf2xpe9ajo1vdw-8rq11"
$ws.Cells.Item(3,5).Value = "Wildlife "
$ws.Cells.Item(3,6).Value = "Monitoring project "
$ws.Cells.Item(3,7).Value = "Ongoing "
$ws.Cells.Item(3,8).Value = "Orange "
$ws.Cells.Item(3,9).Value = "Circle "
$ws.Cells.Item(3,10).Value = "Fill "

# Row 4
$ws.Cells.Item(4,1).Value = "Fish & Crustacean Surveys "
$ws.Cells.Item(4,2).Value = "Fish_Survey_Sites"
$ws.Cells.Item(4,3).Value = "Fish_Survey_Sites"
$ws.Cells.Item(4,4).Value = "jmkhoch.cjsg8g4sn0btk2wntepy64o3k-29zjm"
$ws.Cells.Item(4,5).Value = "Wildlife "
$ws.Cells.Item(4,6).Value = "Monitoring project "
$ws.Cells.Item(4,7).Value = "Ongoing "
$ws.Cells.Item(4,8).Value = "Orange "
$ws.Cells.Item(4,9).Value = "Circle "
$ws.Cells.Item(4,10).Value = "Fill "

# Row 5
$ws.Cells.Item(5,1).Value = "Salt Marsh Plant Survey "
$ws.Cells.Item(5,2).Value = "Vegetation_Survey_Plots_Bronx_Ki"
$ws.Cells.Item(5,3).Value = "Vegetation_Survey_Plots_Bronx_Kill"
$ws.Cells.Item(5,4).Value = "jmkhoch.cjsg8f8xg0xpi2xl66fk8aogz-1q7ip"
$ws.Cells.Item(5,5).Value = "Vegetation "
$ws.Cells.Item(5,6).Value = "Monitoring project "
$ws.Cells.Item(5,7).Value = "Ongoing "
$ws.Cells.Item(5,8).Value = "Green"
$ws.Cells.Item(5,9).Value = "Circle "
$ws.Cells.Item(5,10).Value = "Fill "

# Row 6
$ws.Cells.Item(6,1).Value = "Billion Oyster Project Oyster Cage Monitoring"
$ws.Cells.Item(6,2).Value = "Oyster_Cages"
$ws.Cells.Item(6,3).Value = "Oyster_Cages"
$ws.Cells.Item(6,4).Value = "jmkhoch.cjsg8eq7m116c2xmxg2riqgqr-6x6tv"
$ws.Cells.Item(6,5).Value = "Wildlife "
$ws.Cells.Item(6,6).Value = "Monitoring project "
$ws.Cells.Item(6,7).Value = "Ongoing "
$ws.Cells.Item(6,8).Value = "Orange "
$ws.Cells.Item(6,9).Value = "Circle "
$ws.Cells.Item(6,10).Value = "Fill "

# Row 7
$ws.Cells.Item(7,1).Value = "Salt Marsh Nutrient Cycling Study "
$ws.Cells.Item(7,2).Value = "NutrientCyclingStudy_SampleSites"
$ws.Cells.Item(7,3).Value = "NutrentCyclingStudy_SampleSites "
$ws.Cells.Item(7,4).Value = "jmkhoch.cjsg8btz70f4i2wpglvzxrbt5-1m0jn"
$ws.Cells.Item(7,5).Value = "Water quality "
$ws.Cells.Item(7,6).Value = "Monitoring project "
$ws.Cells.Item(7,7).Value = "Ongoing "
$ws.Cells.Item(7,8).Value = "Blue"
$ws.Cells.Item(7,9).Value = "Circle "
$ws.Cells.Item(7,10).Value = "Fill "

# Row 8
$ws.Cells.Item(8,1).Value = "Citizens' Water Quality Testing Program "
$ws.Cells.Item(8,2).Value = "CWQT_pilot"
$ws.Cells.Item(8,3).Value = "CWQT_pilot"
$ws.Cells.Item(8,4).Value = "jmkhoch.cjr88d2gi056632qcm1k0sypj-3i9cy"
$ws.Cells.Item(8,5).Value = "Water quality "
$ws.Cells.Item(8,6).Value = "Monitoring project "
$ws.Cells.Item(8,7).Value = "Ongoing "
$ws.Cells.Item(8,8).Value = "Blue"
$ws.Cells.Item(8,9).Value = "Circle "
$ws.Cells.Item(8,10).Value = "Fill "

# Row 9
$ws.Cells.Item(9,1).Value = "Oyster Substrate Survey "
$ws.Cells.Item(9,2).Value = "Oyster-Substrate"
$ws.Cells.Item(9,3).Value = "Oyster-Substrate"
$ws.Cells.Item(9,4).Value = "jmkhoch.cjsvzh59w0gik32rsmgau3drp-60ebk"
$ws.Cells.Item(9,5).Value = "Wildlife "
$ws.Cells.Item(9,6).Value = "Academic partner project "
$ws.Cells.Item(9,7).Value = "Completed "
$ws.Cells.Item(9,8).Value = "Orange "
$ws.Cells.Item(9,9).Value = "Triangle"
$ws.Cells.Item(9,10).Value = "Hollow Fill "

# Row 10
$ws.Cells.Item(10,1).Value = "Fish & Water Quality Survey "
$ws.Cells.Item(10,2).Value = "Fish-WQ-Project-2019"
$ws.Cells.Item(10,3).Value = "Fish-WQ-Project-2019"
$ws.Cells.Item(10,4).Value = "jmkhoch.cjsvzoelh06nh2wpdify23aus-3k9fj"
$ws.Cells.Item(10,5).Value = "Wildlife "
$ws.Cells.Item(10,6).Value = "Academic partner project "
$ws.Cells.Item(10,7).Value = "Completed "
$ws.Cells.Item(10,8).Value = "Orange "
$ws.Cells.Item(10,9).Value = "Triangle"
$ws.Cells.Item(10,10).Value = "Hollow Fill "

# Row 11
$ws.Cells.Item(11,1).Value = "Wild Oyster Population Survey "
$ws.Cells.Item(11,2).Value = "Wild-Oyster-Survey"
$ws.Cells.Item(11,3).Value = "Wild-Oyster-Survey"
$ws.Cells.Item(11,4).Value = "jmkhoch.cjsvzo5zq0iy82wrsjor40rkd-1l6yt"
$ws.Cells.Item(11,5).Value = "Wildlife "
$ws.Cells.Item(11,6).Value = "Academic partner project "
$ws.Cells.Item(11,7).Value = "Completed "
$ws.Cells.Item(11,8).Value = "Orange "
$ws.Cells.Item(11,9).Value = "Triangle"
$ws.Cells.Item(11,10).Value = "Hollow Fill "

# Row 12
$ws.Cells.Item(12,1).Value = "Benthic MacroInvertebrate Survey of the Little Hell Gate Salt Marsh "
$ws.Cells.Item(12,2).Value = "Benthic-Macroinvertebrate-Survey"
$ws.Cells.Item(12,3).Value = "Benthic-Macroinvertebrate-Survey"
$ws.Cells.Item(12,4).Value = "jmkhoch.cjsvznykb02t22so2eomqip06-24ueg"
$ws.Cells.Item(12,5).Value = "Wildlife "
$ws.Cells.Item(12,6).Value = "Academic partner project "
$ws.Cells.Item(12,7).Value = "Completed "
$ws.Cells.Item(12,8).Value = "Orange "
$ws.Cells.Item(12,9).Value = "Triangle"
$ws.Cells.Item(12,10).Value = "Hollow Fill "

# Row 13
$ws.Cells.Item(13,1).Value = "Pollinator Surveys "
$ws.Cells.Item(13,2).Value = "pollinator-survey"
$ws.Cells.Item(13,3).Value = "pollinator-survey"
$ws.Cells.Item(13,4).Value = "jmkhoch.cjsvzgcmm0dhi2rqsls6msy0d-4f4n2"
$ws.Cells.Item(13,5).Value = "Wildlife "
$ws.Cells.Item(13,6).Value = "Monitoring project "
$ws.Cells.Item(13,7).Value = "Ongoing "
$ws.Cells.Item(13,8).Value = "Orange "
$ws.Cells.Item(13,9).Value = "Circle "
$ws.Cells.Item(13,10).Value = "Fill "

# --- Restore view state: active selection at G5 (scroll position B1 not exposed via this COM shim) ---
$ws.Range("G5").Select()
